# Auto-generated script applying cell value corrections per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 77412.75
$ws.Range("J3").Value = 77412.75
$ws.Range("L3").Value = 77412.75
$ws.Range("N3").Value = -77640.75
$ws.Range("H8").Value = 7.5
$ws.Range("I8").Value = 7.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 22.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 116.5
$ws.Range("N8").ClearContents()
$ws.Range("H33").Value = 1009.36365
$ws.Range("J33").Value = 1169.5
$ws.Range("L33").Value = 1169.5
$ws.Range("N33").Value = -1627.5
$ws.Range("H82").Value = 6354.875
$ws.Range("I82").Value = 5119.857
$ws.Range("J82").Value = 15000
$ws.Range("K82").Value = 15359.571
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -14953.571
$ws.Range("N82").Value = -45812
$ws.Range("H85").Value = 6354.875
$ws.Range("I85").Value = 5119.857
$ws.Range("J85").Value = 15000
$ws.Range("K85").Value = 15359.571
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -13955.571
$ws.Range("N85").Value = -47808
$ws.Range("H86").Value = 4053330
$ws.Range("I86").Value = 4773.375
$ws.Range("K86").Value = 4773.375
$ws.Range("M86").Value = -3650.375
$ws.Range("H89").Value = 4053330
$ws.Range("I89").Value = 4773.375
$ws.Range("K89").Value = 23866.875
$ws.Range("M89").Value = -18250.875
$ws.Range("H100").Value = 4673.3335
$ws.Range("J100").Value = 7833.3335
$ws.Range("L100").Value = 7833.3335
$ws.Range("N100").Value = -8915.333500000001
$ws.Range("H102").Value = 77412.75
$ws.Range("J102").Value = 77412.75
$ws.Range("L102").Value = 77412.75
$ws.Range("N102").Value = -83902.75
$ws.Range("H112").Value = 1586.9818
$ws.Range("J112").Value = 1589.0555
$ws.Range("L112").Value = 4767.166499999999
$ws.Range("N112").Value = -6983.166499999999
$ws.Range("H113").Value = 8342.25
$ws.Range("I113").Value = 6147.6
$ws.Range("K113").Value = 6147.6
$ws.Range("M113").Value = -2893.6
$ws.Range("H115").Value = 903.7059
$ws.Range("I115").Value = 827.5333000000001
$ws.Range("J115").Value = 1475
$ws.Range("K115").Value = 2482.5999
$ws.Range("L115").Value = 4425
$ws.Range("M115").Value = -915.5999000000002
$ws.Range("N115").Value = -7559
$ws.Range("H135").Value = 71434720
$ws.Range("I135").Value = 71434720
$ws.Range("K135").Value = 642912480
$ws.Range("M135").Value = -642909945
$ws.Range("H137").Value = 1953
$ws.Range("I137").Value = 1363.9
$ws.Range("J137").Value = 3916.6667
$ws.Range("K137").Value = 4091.7
$ws.Range("L137").Value = 11750.0001
$ws.Range("M137").Value = -1541.7
$ws.Range("N137").Value = -16850.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1201
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 3400.1
$ws.Range("K61").Value = 3400.1
$ws.Range("M61").Value = -3188.1
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 3400.1
$ws.Range("K136").Value = 10200.3
$ws.Range("M136").Value = -7650.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H105").Value = 6501.8
$ws.Range("I105").Value = 4254.5
$ws.Range("K105").Value = 4254.5
$ws.Range("M105").Value = -2507.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4313.048
$ws.Range("I31").Value = 1954.7778
$ws.Range("J31").Value = 6081.75
$ws.Range("K31").Value = 1954.7778
$ws.Range("L31").Value = 6081.75
$ws.Range("M31").Value = -1659.7778
$ws.Range("N31").Value = -6671.75
$ws.Range("H34").Value = 4313.048
$ws.Range("I34").Value = 1954.7778
$ws.Range("J34").Value = 6081.75
$ws.Range("K34").Value = 1954.7778
$ws.Range("L34").Value = 6081.75
$ws.Range("M34").Value = -1752.7778
$ws.Range("N34").Value = -6485.75
$ws.Range("H58").Value = 3312.2666
$ws.Range("J58").Value = 5557.6
$ws.Range("L58").Value = 5557.6
$ws.Range("N58").Value = -5963.6
$ws.Range("H99").Value = 10271.4
$ws.Range("J99").Value = 9746
$ws.Range("L99").Value = 9746
$ws.Range("N99").Value = -12742
$ws.Range("H126").Value = 10271.4
$ws.Range("J126").Value = 9746
$ws.Range("L126").Value = 29238
$ws.Range("N126").Value = -34178
$ws.Range("H132").Value = 2507.2
$ws.Range("I132").Value = 2134.6875
$ws.Range("K132").Value = 6404.0625
$ws.Range("M132").Value = -3874.0625
$ws.Range("H134").Value = 3765.54
$ws.Range("I134").Value = 2613.7058
$ws.Range("K134").Value = 7841.117400000001
$ws.Range("M134").Value = -5306.117400000001
$ws.Range("H136").Value = 3312.2666
$ws.Range("J136").Value = 5557.6
$ws.Range("L136").Value = 16672.8
$ws.Range("N136").Value = -21772.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58.526318
$ws.Range("I2").Value = 56.5
$ws.Range("K2").Value = 339
$ws.Range("M2").Value = -226
$ws.Range("H17").Value = 2927
$ws.Range("I17").Value = 2927
$ws.Range("K17").Value = 8781
$ws.Range("M17").Value = -8612
$ws.Range("H63").Value = 173033.33
$ws.Range("I63").Value = 1002488
$ws.Range("K63").Value = 3007464
$ws.Range("M63").Value = -3006715
$ws.Range("H66").Value = 173033.33
$ws.Range("I66").Value = 1002488
$ws.Range("K66").Value = 9022392
$ws.Range("M66").Value = -9018648
$ws.Range("H130").Value = 502985.44
$ws.Range("I130").Value = 876474.5
$ws.Range("K130").Value = 2629423.5
$ws.Range("M130").Value = -2624403.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1112559.9
$ws.Range("I3").Value = 522.5
$ws.Range("J3").Value = 3336634.8
$ws.Range("K3").Value = 522.5
$ws.Range("L3").Value = 3336634.8
$ws.Range("M3").Value = -406.5
$ws.Range("N3").Value = -3336866.8
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H104").Value = 59999.5
$ws.Range("J104").Value = 59999.5
$ws.Range("L104").Value = 59999.5
$ws.Range("N104").Value = -66987.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7158.8335
$ws.Range("I7").Value = 9438
$ws.Range("J7").Value = 5335.5
$ws.Range("K7").Value = 9438
$ws.Range("L7").Value = 5335.5
$ws.Range("M7").Value = -9326
$ws.Range("N7").Value = -5559.5
$ws.Range("H40").Value = 5892.4443
$ws.Range("I40").Value = 2621.4
$ws.Range("J40").Value = 7150.5386
$ws.Range("K40").Value = 2621.4
$ws.Range("L40").Value = 7150.5386
$ws.Range("M40").Value = -2485.4
$ws.Range("N40").Value = -7422.5386
$ws.Range("H126").Value = 7158.8335
$ws.Range("I126").Value = 9438
$ws.Range("J126").Value = 5335.5
$ws.Range("K126").Value = 28314
$ws.Range("L126").Value = 16006.5
$ws.Range("M126").Value = -25844
$ws.Range("N126").Value = -20946.5
$ws.Range("H132").Value = 4267.4443
$ws.Range("I132").Value = 3432.5881
$ws.Range("K132").Value = 10297.7643
$ws.Range("M132").Value = -7767.764299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 53000
$ws.Range("J11").Value = 97000
$ws.Range("L11").Value = 97000
$ws.Range("N11").Value = -97284
$ws.Range("H33").Value = 11438.333
$ws.Range("I33").Value = 2410
$ws.Range("J33").Value = 29495
$ws.Range("K33").Value = 2410
$ws.Range("L33").Value = 29495
$ws.Range("M33").Value = -2160
$ws.Range("N33").Value = -29995
$ws.Range("H36").Value = 11438.333
$ws.Range("I36").Value = 2410
$ws.Range("J36").Value = 29495
$ws.Range("K36").Value = 2410
$ws.Range("L36").Value = 29495
$ws.Range("M36").Value = -2160
$ws.Range("N36").Value = -29995
$ws.Range("H96").Value = 2102.9092
$ws.Range("J96").Value = 2200
$ws.Range("L96").Value = 2200
$ws.Range("N96").Value = -4946
$ws.Range("H107").Value = 979.2
$ws.Range("I107").Value = 979.2
$ws.Range("K107").Value = 2937.6
$ws.Range("M107").Value = -1017.6
$ws.Range("H119").Value = 87000
$ws.Range("J119").Value = 87000
$ws.Range("L119").Value = 87000
$ws.Range("N119").Value = -96676
$ws.Range("H132").Value = 2617.4211
$ws.Range("J132").Value = 3175.7778
$ws.Range("L132").Value = 9527.3334
$ws.Range("N132").Value = -14587.3334
